$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers must now be upper case
$ws.Range("A2").Value = "STATION"
$ws.Range("B2").Value = "NAME"
$ws.Range("C2").Value = "NETID"
$ws.Range("D2").Value = "LON"
$ws.Range("E2").Value = "INTENSITY"

# Update selected cell to E2
$ws.Range("E2").Select()
